$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'51.128.22"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.33%  '

$ws.Range('D3').Value = "'2.963.43"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.98%  '

$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').Value = "'381.08"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.34%  '

$ws.Range('D6').Value = "'102.65"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.45%  '

$ws.Range('D7').Value = "'0.545"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.77%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('D9').Value = "'0.589"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.26%  '

$ws.Range('D10').Value = "'36.56"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.08%  '

$ws.Range('E11').Value = '  -0.35%  '

$ws.Range('D12').Value = "'0.0853"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.36%  '

$ws.Range('B13').Value = 'Uniswap'
$ws.Range('C13').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D13').Value = "'12.51"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +75.90%  '

$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = "'18.47"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.16%  '

$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = "'3.418.05"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.65%  '

$ws.Range('D16').Value = "'7.77"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +6.03%  '

$ws.Range('D17').Value = "'2.961.16"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.97%  '

$ws.Range('D18').Value = "'1.02"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.81%  '

$ws.Range('D19').Value = "'51.163.89"
$ws.Range('D19').Style = 'Normal'

$ws.Range('D20').Value = "'3.07"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.83%  '

$ws.Range('D21').Value = "'12.42"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.22%  '

$ws.Range('D22').Value = "'0.0₃0958"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.38%  '

$ws.Range('D23').Value = "'3.37"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +18.12%  '

$ws.Range('E24').Value = '  +2.94%  '

$ws.Range('E25').Value = '  +2.37%  '

$ws.Range('D26').Value = "'8.02"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.58%  '

$ws.Range('E27').Value = '  +0.09%  '

$ws.Range('D28').Value = "'0.167"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.03%  '

$ws.Range('D29').Value = "'25.94"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.41%  '

$ws.Range('D30').Value = "'7.01"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -9.13%  '

$ws.Range('E31').Value = '  -3.80%  '

$ws.Range('D32').Value = "'10.57"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.68%  '

$ws.Range('B33').Value = 'OKB'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D33').Value = "'51.23"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.20%  '

$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D34').Value = "'34.22"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.86%  '

$ws.Range('D35').Value = "'2.06"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.30%  '

$ws.Range('D36').Value = "'0.0437"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.10%  '

$ws.Range('E37').Value = '  +0.01%  '

$ws.Range('D38').Value = "'3.25"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +9.61%  '

$ws.Range('D39').Value = "'16.79"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.69%  '

$ws.Range('E40').Value = '  +2.31%  '

$ws.Range('D41').Value = "'1.84"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.38%  '

$ws.Range('E42').Value = '  -1.92%  '

$ws.Range('D43').Value = "'124.55"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.72%  '

$ws.Range('D44').Value = "'3.61"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +12.33%  '

$ws.Range('D45').Value = "'21.84"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.76%  '

$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = "'2.102.44"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.91%  '

$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').Value = "'2.03"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.81%  '

$ws.Range('E48').Value = '  -0.83%  '

$ws.Range('D49').Value = "'0.259"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.96%  '

$ws.Range('D50').Value = "'0.0324"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.98%  '

$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').Value = "'1.34"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +7.49%  '

